$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","K","M","N","O")

$data = @{
    2 = @(0.0494028206232997,0.01849215619530753,0.107369397667668,0.6985174721780538,0.002411493739661308,0.9142017256984332,0.3447736560447296,1.103131043754892,2.354599466279382)
    3 = @(0.04386124092449961,0.01751524613989019,0.1015636426510795,0.6901919500866853,0.002414581642825512,0.8027137754717444,0.3090260743027002,1.114037838243455,2.341576173211934)
    4 = @(0.0404731889614709,0.01691130418702258,0.09809643650871891,0.6856156552827741,0.002416578722254152,0.7341795315847719,0.2871659319520035,1.121206729647916,2.335385109436857)
    5 = @(0.03909615488895213,0.0166641772647047,0.09670779273211139,0.6838851389123519,0.002417418048042136,0.7062322367071943,0.2782800785363335,1.124246844273721,2.333314982129366)
    6 = @(0.03886771772241104,0.01662308124809186,0.09647866783028292,0.6836058936459395,0.002417558959779117,0.7015904935571484,0.2768059367284508,1.124758825056389,2.332998551684881)
    7 = @(0.04045460314949878,0.0169079754340693,0.09807761084293531,0.685591773340704,0.002416589938364428,0.7338027001553371,0.2870460038690297,1.121247248916205,2.335355359244261)
    8 = @(0.0474890597087807,0.01815618287676202,0.1053471686378629,0.6955354675602905,0.002412537513711639,0.8757779669547858,0.332429371143661,1.106793741737818,2.349733628482909)
    9 = @(0.06140034188047139,0.02057053502860384,0.1203892729417575,0.7193015153609537,0.002405389270951043,1.153518834199531,0.4221394474421913,1.0821951469472,2.392311774403794)
    10 = @(0.07169546622115774,0.02232320280938183,0.1319386693923832,0.7393895485294024,0.002400619240189918,1.35713731751423,0.4885010977999684,1.066404232780982,2.432450070841099)
    11 = @(0.07639593772768194,0.02311578368220779,0.1373049757959777,0.7491043898836693,0.002398552772291443,1.449669001012012,0.5187928424108321,1.059715861105538,2.452652551826418)
    12 = @(0.0781783915498977,0.02341521841128724,0.139353528304639,0.7528664842072033,0.002397785048334985,1.484693784582419,0.5302785973174622,1.057254342187044,2.460583594538718)
    13 = @(0.07779439731173454,0.02335076116869317,0.1389115998071091,0.7520525394815252,0.002397949734164428,1.477151255647641,0.5278042697825356,1.057781305952389,2.458862992370342)
    14 = @(0.07654253132216127,0.02314043247028508,0.1374731801681008,0.7494122281951832,0.002398489314873251,1.452550817871895,0.5197374829166819,1.05951192279597,2.453299408073519)
    15 = @(0.07577605172095048,0.02301150851317857,0.1365942571002634,0.7478058196200408,0.002398821749170676,1.437480354904608,0.5147982885868032,1.060581252197913,2.449928160742303)
    16 = @(0.07138862640597665,0.02227130914559439,0.1315902550416226,0.7387662967125408,0.002400756363951905,1.351088128781839,0.4865235475484724,1.066851297660278,2.43116901187804)
    17 = @(0.06870149344855747,0.02181599836249148,0.1285494574242563,0.73336880835285,0.002401969630200752,1.2980639808913,0.4692044497779193,1.070824575821881,2.420159583034007)
    18 = @(0.06715754046557265,0.02155367290127685,0.126811036403943,0.7303185750534738,0.002402677209725571,1.267556988999388,0.4592527015725665,1.073156488585667,2.41401008905757)
    19 = @(0.06663506129076779,0.0214647785854396,0.1262242429599567,0.7292951263777638,0.002402918459101885,1.257226354287411,0.4558848819007011,1.073954034562483,2.411959337709618)
    20 = @(0.06898737607717464,0.02186451294115699,0.1288720600129238,0.7339377619163514,0.00240183946844162,1.303709422541374,0.4710470890538829,1.07039679157549,2.421312622687537)
    21 = @(0.07691016705808806,0.02320223022314849,0.137895230218362,0.7501854881018915,0.002398330425882551,1.459776978543175,0.5221064882255035,1.059001665355105,2.454925936016963)
    22 = @(0.08210267518074943,0.02407242546159694,0.1438884010578363,0.7612899923584422,0.002396123312226823,1.561688398218678,0.5555637956273216,1.051969453745059,2.478531533333296)
    23 = @(0.07933000247155064,0.02360836587409665,0.1406808548732315,0.7553187487980182,0.002397293422319558,1.507304792717662,0.5376990212422754,1.055684674376508,2.465782517881195)
    24 = @(0.06885812569213101,0.02184258124636784,0.1287261809611238,0.7336803733965098,0.002401898283251938,1.301157189000151,0.4702140155498427,1.070590044488966,2.420790772942951)
    25 = @(0.05762412821110274,0.01992104640781633,0.1162338769982725,0.7124127108677669,0.002407238092294347,1.078457449965811,0.3977923958518588,1.088448971057915,2.379244413025418)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
